$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the keyword/appID pairs that changed content (new shared-string
# order: "powerful quotes", "motivation quotes", "optimistic quotes").
$ws.Range("A8").Value = "powerful quotes"
$ws.Range("A12").Value = "powerful quotes"

$ws.Range("A15").Value = "motivation quotes"
$ws.Range("A18").Value = "motivation quotes"
$ws.Range("B18").Value = "com.sugar.powerfulquotes"
$ws.Range("A20").Value = "motivation quotes"
$ws.Range("B20").Value = "com.sugar.powerfulquotes"
$ws.Range("A21").Value = "motivation quotes"
$ws.Range("B21").Value = "com.sugar.powerfulquotes"

$ws.Range("A9").Value = "optimistic quotes"
$ws.Range("B9").Value = "com.sugar.powerfulquotes"
$ws.Range("A14").Value = "optimistic quotes"
$ws.Range("B14").Value = "com.sugar.powerfulquotes"

# Update the active selection to match the author's final cursor position.
$ws.Range("B21").Select()
